$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repull data / push all data / mean calculation: update dSF (column F) values
$ws.Range("F3").Value = 5
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = 5
$ws.Range("F9").Value = 1
$ws.Range("F11").Value = 3
$ws.Range("F14").Value = -4
$ws.Range("F18").Value = 0
